# Scheduled-runner data refresh: updates market-price-derived columns
# (H..N) for a handful of Leve rows across the ALC/ARM/BSM/CRP/CUL/LTW/WVR
# sheets. Values come from an external pricing source, not in-sheet
# formulas, so each target cell is written explicitly.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 398.6842
$ws.Range("J17").Value = 404.05405
$ws.Range("L17").Value = 1212.16215
$ws.Range("N17").Value = -1548.16215
$ws.Range("H113").Value = 4168679.5
$ws.Range("I113").Value = 5884100.5
$ws.Range("J113").Value = 2657
$ws.Range("K113").Value = 5884100.5
$ws.Range("L113").Value = 2657
$ws.Range("M113").Value = -5880846.5
$ws.Range("N113").Value = -9165
$ws.Range("H127").Value = 1879
$ws.Range("I127").Value = 480.875
$ws.Range("J127").Value = 2811.0833
$ws.Range("K127").Value = 1442.625
$ws.Range("L127").Value = 8433.249899999999
$ws.Range("M127").Value = 3517.375
$ws.Range("N127").Value = -18353.2499
$ws.Range("H132").Value = 2058845.1
$ws.Range("I132").Value = 1255.2264
$ws.Range("J132").Value = 111111110
$ws.Range("K132").Value = 3765.6792
$ws.Range("L132").Value = 333333330
$ws.Range("M132").Value = -1235.6792
$ws.Range("N132").Value = -333338390
$ws.Range("H137").Value = 17375640
$ws.Range("I137").Value = 988.1786
$ws.Range("J137").Value = 78186930
$ws.Range("K137").Value = 2964.5358
$ws.Range("L137").Value = 234560790
$ws.Range("M137").Value = -414.5357999999997
$ws.Range("N137").Value = -234565890
$ws.Range("H138").Value = 2539.5542
$ws.Range("I138").Value = 1768
$ws.Range("K138").Value = 5304
$ws.Range("M138").Value = -164
$ws.Range("H141").Value = 1625.1852
$ws.Range("I141").Value = 1278.3334
$ws.Range("J141").Value = 4400
$ws.Range("K141").Value = 3835.0002
$ws.Range("L141").Value = 13200
$ws.Range("M141").Value = 1344.9998
$ws.Range("N141").Value = -23560
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6031545
$ws.Range("I32").Value = 7778.239
$ws.Range("K32").Value = 7778.239
$ws.Range("M32").Value = -7491.239
$ws.Range("H74").Value = 32625500
$ws.Range("I74").Value = 29412452
$ws.Range("J74").Value = 41028856
$ws.Range("K74").Value = 29412452
$ws.Range("L74").Value = 41028856
$ws.Range("M74").Value = -29411578
$ws.Range("N74").Value = -41030604
$ws.Range("H77").Value = 32625500
$ws.Range("I77").Value = 29412452
$ws.Range("J77").Value = 41028856
$ws.Range("K77").Value = 147062260
$ws.Range("L77").Value = 205144280
$ws.Range("M77").Value = -147057892
$ws.Range("N77").Value = -205153016
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1994.05
$ws.Range("I86").Value = 1994
$ws.Range("J86").Value = 1999
$ws.Range("K86").Value = 1994
$ws.Range("L86").Value = 1999
$ws.Range("M86").Value = -871
$ws.Range("N86").Value = -4245
$ws.Range("H89").Value = 1994.05
$ws.Range("I89").Value = 1994
$ws.Range("J89").Value = 1999
$ws.Range("K89").Value = 9970
$ws.Range("L89").Value = 9995
$ws.Range("M89").Value = -4354
$ws.Range("N89").Value = -21227
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1605676.5
$ws.Range("I31").Value = 980.25
$ws.Range("J31").Value = 8941431
$ws.Range("K31").Value = 980.25
$ws.Range("L31").Value = 8941431
$ws.Range("M31").Value = -685.25
$ws.Range("N31").Value = -8942021
$ws.Range("H34").Value = 1605676.5
$ws.Range("I34").Value = 980.25
$ws.Range("J34").Value = 8941431
$ws.Range("K34").Value = 980.25
$ws.Range("L34").Value = 8941431
$ws.Range("M34").Value = -778.25
$ws.Range("N34").Value = -8941835
$ws.Range("H134").Value = 635997.4399999999
$ws.Range("I134").Value = 973.6923
$ws.Range("J134").Value = 3637928
$ws.Range("K134").Value = 2921.0769
$ws.Range("L134").Value = 10913784
$ws.Range("M134").Value = -386.0769
$ws.Range("N134").Value = -10918854
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2506702.5
$ws.Range("I5").Value = 3344967.2
$ws.Range("J5").Value = 1588603.4
$ws.Range("K5").Value = 10034901.6
$ws.Range("L5").Value = 4765810.199999999
$ws.Range("M5").Value = -10034789.6
$ws.Range("N5").Value = -4766034.199999999
$ws.Range("H8").Value = 143.5238
$ws.Range("I8").Value = 143.5238
$ws.Range("K8").Value = 430.5714
$ws.Range("M8").Value = -291.5714
$ws.Range("H32").Value = 1789.0769
$ws.Range("I32").Value = 1451.6
$ws.Range("J32").Value = 2000
$ws.Range("K32").Value = 4354.799999999999
$ws.Range("L32").Value = 6000
$ws.Range("M32").Value = -4071.799999999999
$ws.Range("N32").Value = -6566
$ws.Range("H112").Value = 3125
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 3125
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 9375
$ws.Range("M112").ClearContents()
$ws.Range("N112").Value = -11591
$ws.Range("H122").Value = 1719.7858
$ws.Range("I122").Value = 302.83334
$ws.Range("K122").Value = 2725.50006
$ws.Range("M122").Value = -275.5000600000003
$ws.Range("H131").Value = 8334121.5
$ws.Range("I131").Value = 31250354
$ws.Range("J131").Value = 946.4545000000001
$ws.Range("K131").Value = 93751062
$ws.Range("L131").Value = 2839.3635
$ws.Range("M131").Value = -93746022
$ws.Range("N131").Value = -12919.3635
$ws.Range("H135").Value = 2506702.5
$ws.Range("I135").Value = 3344967.2
$ws.Range("J135").Value = 1588603.4
$ws.Range("K135").Value = 30104704.8
$ws.Range("L135").Value = 14297430.6
$ws.Range("M135").Value = -30102169.8
$ws.Range("N135").Value = -14302500.6
$ws.Range("H137").Value = 2855.4443
$ws.Range("I137").Value = 3412
$ws.Range("J137").Value = 2641.3845
$ws.Range("K137").Value = 10236
$ws.Range("L137").Value = 7924.1535
$ws.Range("M137").Value = -5136
$ws.Range("N137").Value = -18124.1535
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 1985054.8
$ws.Range("I136").Value = 2137689.8
$ws.Range("J136").Value = 798.75
$ws.Range("K136").Value = 6413069.399999999
$ws.Range("L136").Value = 2396.25
$ws.Range("M136").Value = -6410519.399999999
$ws.Range("N136").Value = -7496.25
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 702047.8
$ws.Range("I132").Value = 1813.3684
$ws.Range("J132").Value = 4330535.5
$ws.Range("K132").Value = 5440.1052
$ws.Range("L132").Value = 12991606.5
$ws.Range("M132").Value = -2910.1052
$ws.Range("N132").Value = -12996666.5
$ws.Range("H136").Value = 820.54877
$ws.Range("I136").Value = 304.33963
$ws.Range("J136").Value = 1763.9656
$ws.Range("K136").Value = 913.0188900000001
$ws.Range("L136").Value = 5291.8968
$ws.Range("M136").Value = 1636.98111
$ws.Range("N136").Value = -10391.8968
